$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: Soundtrack comments ---
$ws.Range("F47").Value = 5
$ws.Range("G47").Value = "Soundtrack is the medley of Galaga from Smash Brothers Ultimate (in Galaga)"
$ws.Range("H47").Value = "A few audio effects from galaga included"

# --- Row 48: audio effects detail ---
$ws.Range("F48").Value = 5
$ws.Range("G48").Value = "shooting, enemy descent, and bullets killing enemies"

# --- Row 17: scaling transformation ---
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = "All assets being used are automatically scaled to better fit the screen"

# --- Row 18: shearing transformation ---
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = "The bullets are being sheared to be like the game (straight bullets)"

# --- Row 12: broad phase collision detection (BSP tree) ---
$ws.Range("F12").Value = 40
$ws.Range("G12").Value = "Implemented broad phase collision detection using BSP tree (or at least what I imagine it to be for galaga)"
$ws.Range("H12").Value = "Essentially, 3 lists of bools getting down to very specific numbers of where gameobjects are"

# --- Update active cell selection ---
[void]$ws.Range("K4").Select()
